$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from E1 (header cell) onto the new F1 header cell
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Minor floating point updates to existing values
$ws.Range("B2").Value = 0.07514644587374561
$ws.Range("D2").Value = 0.2119198634755611

# New model name cell
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"
